$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'" + $text
    $c.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 "244.24"
Set-TextValue 2 5 "-1.21%"

# Row 3
Set-TextValue 3 4 "27.38"
Set-TextValue 3 5 "3.84%"

# Row 4
Set-TextValue 4 5 "-0.59%"

# Row 5
Set-TextValue 5 4 "0.05682"
Set-TextValue 5 5 "1.07%"

# Row 6
Set-TextValue 6 4 "6.475"
Set-TextValue 6 5 "-0.62%"

# Row 7
Set-TextValue 7 4 "0.8214"
Set-TextValue 7 5 "0.81%"

# Row 8
Set-TextValue 8 4 "0.8406"
Set-TextValue 8 5 "-0.97%"

# Row 9
Set-TextValue 9 2 "WazirX"
Set-TextValue 9 3 "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue 9 4 "0.1325"
Set-TextValue 9 5 "-1.65%"

# Row 10
Set-TextValue 10 2 "MandalaExchangeToken"
Set-TextValue 10 3 "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue 10 4 "0.06913"
Set-TextValue 10 5 "-0.78%"

# Row 11
Set-TextValue 11 2 "BitrueCoin"
Set-TextValue 11 3 "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue 11 4 "0.02860"
Set-TextValue 11 5 "1.37%"

# Row 12
Set-TextValue 12 2 "BitMartToken"
Set-TextValue 12 3 "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue 12 4 "0.09398"
Set-TextValue 12 5 "-0.11%"

# Row 13
Set-TextValue 13 2 "BitForexToken"
Set-TextValue 13 3 "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue 13 4 "0.001511"
Set-TextValue 13 5 "0.03%"

# Row 14
Set-TextValue 14 2 "CoinExToken"
Set-TextValue 14 3 "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue 14 4 "0.04106"
Set-TextValue 14 5 "-11.90%"

# Row 15
Set-TextValue 15 2 "One"
Set-TextValue 15 3 "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue 15 4 "0.0006019"
Set-TextValue 15 5 "0.98%"

# Row 16
Set-TextValue 16 4 "0.006139"
Set-TextValue 16 5 "0.17%"

# Row 17
Set-TextValue 17 5 "-2.21%"

# Row 18
Set-TextValue 18 4 "3.000"
Set-TextValue 18 5 "-1.93%"

# Row 19
Set-TextValue 19 4 "2.224"
Set-TextValue 19 5 "4.98%"

# Row 20
Set-TextValue 20 4 "0.3152"
Set-TextValue 20 5 "-0.91%"

# Row 21
Set-TextValue 21 4 "0.03194"
Set-TextValue 21 5 "-0.43%"

# Row 22
Set-TextValue 22 4 "0.1291"
Set-TextValue 22 5 "-2.19%"

# Row 23
Set-TextValue 23 4 "3.561"
Set-TextValue 23 5 "-5.19%"

# Row 25
Set-TextValue 25 4 "0.001219"
Set-TextValue 25 5 "-2.32%"

# Row 26
Set-TextValue 26 4 "0.003947"
Set-TextValue 26 5 "-14.43%"

# Row 27
Set-TextValue 27 4 "0.00009800"
Set-TextValue 27 5 "2.09%"

# Row 40
Set-TextValue 40 4 "0.03789"
Set-TextValue 40 5 "3.12%"

# Row 41
Set-TextValue 41 4 "0.006180"
Set-TextValue 41 5 "-0.64%"

# Row 42
Set-TextValue 42 4 "0.1054"
Set-TextValue 42 5 "-0.59%"

# Row 43
Set-TextValue 43 4 "0.002410"
Set-TextValue 43 5 "-3.59%"

# Row 44
Set-TextValue 44 4 "0.009689"
Set-TextValue 44 5 "11.78%"

# Row 45
Set-TextValue 45 4 "0.00005209"
Set-TextValue 45 5 "-1.60%"

# Row 46
Set-TextValue 46 5 "-0.01%"

# Row 47
Set-TextValue 47 5 "-15.44%"

# Row 48
Set-TextValue 48 4 "0.002570"
Set-TextValue 48 5 "0.93%"

# Row 49
Set-TextValue 49 5 "-0.01%"

# Row 50
Set-TextValue 50 5 "-0.01%"
